$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja1 (sheet1): currently has header row, a 2002 data row (row 2) and a
# left-over blank-but-styled row (row 7). Fill in the missing seasons
# 2003-2010 in order, clearing the stray style on the old row 7 so it
# becomes a normal text cell like the rest of the list.
$years1 = @(2003,2004,2005,2006,2007,2008,2009,2010)

$ws1.Cells.Item(7,1).Style = "Normal"
$ws1.Cells.Item(7,2).NumberFormat = "@"

for ($i = 0; $i -lt $years1.Length; $i++) {
    $row = $i + 3
    $year = $years1[$i]
    $ws1.Cells.Item($row, 1).Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\$year.xlsx"
    $ws1.Cells.Item($row, 2).NumberFormat = "@"
    $ws1.Cells.Item($row, 2).Value = "$year"
}

# --- Hoja2 (sheet2): currently starts the list at 2003 (row 2) and runs
# through 2008 (row 7, which carries a stray one-off style). Insert the
# missing 2002 season at the top and append 2009/2010 at the bottom, and
# normalise the old "last row" style now that 2008 is a mid-list row.
$ws2.Rows.Item(2).Insert()

$ws2.Cells.Item(2,1).Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2002.xlsx"
$ws2.Cells.Item(2,2).NumberFormat = "@"
$ws2.Cells.Item(2,2).Value = "2002"

$ws2.Cells.Item(8,1).Style = "Normal"
$ws2.Cells.Item(8,1).Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2008.xlsx"
$ws2.Cells.Item(8,2).Style = "Normal"
$ws2.Cells.Item(8,2).NumberFormat = "@"
$ws2.Cells.Item(8,2).Value = "2008"

$ws2.Cells.Item(9,1).Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2009.xlsx"
$ws2.Cells.Item(9,2).NumberFormat = "@"
$ws2.Cells.Item(9,2).Value = "2009"

$ws2.Cells.Item(10,1).Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2010.xlsx"
$ws2.Cells.Item(10,2).NumberFormat = "@"
$ws2.Cells.Item(10,2).Value = "2010"

# --- Selections: both sheets now span A2:B10; keep Hoja1 as the active tab
# (activate Hoja2 first so the final Activate leaves Hoja1 selected).
[void]$ws2.Activate()
[void]$ws2.Range("A2:B10").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A2:B10").Select()
